$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: strip transient <w:proofErr/> markers (and merge the runs they
# split) from the paragraph that contains the given Range, by round-tripping
# the paragraph's Range through WordOpenXML -> InsertXML. The engine's
# InsertXML rebuild does not re-emit proofErr markers, which both removes
# stray gramStart/gramEnd|spellStart/spellEnd wrappers and coalesces runs
# that were only split because of them.
# ---------------------------------------------------------------------------
function CleanParagraphProofErr($para) {
    $prng = $para.Range
    $xml = $prng.WordOpenXML
    [void]$prng.InsertXML($xml)
}

function FindAndClean($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $para = $rng.Paragraphs(1)
    CleanParagraphProofErr $para
}

# 1) Merge the Dynata compensation paragraph (3x spellStart/spellEnd around
#    "Dynata") into a single run. The paragraph's runs carry a leftover
#    w:rsidRPr attribute in the source; a straight WordOpenXML round-trip
#    would keep it, but the target is a bare <w:r>, so do a literal
#    Find/Replace of the whole paragraph text instead - that rebuilds the
#    run (and drops the stale rsid attribute) the way Word does on an
#    in-place retype.
$dynataText = "This survey is made available to respondents via Dynata, which offers great diversity in incentives as some people are motivated by cash, points, or by being able to donate to charity. Others are motivated by the chance to make a difference, make their voice heard, have fun taking a survey, helping out, or having a say in the products and services of the future. Others are motivated by learning opportunities provided by the survey or by the promise of receiving information after taking it. Dynata aims to respond to all of these individual motivations in order to provide a sample which is diverse and as representative as possible of the target population. Dynata uses a reasonable level of reward based on the amount of effort required, the population, and appropriate regional customs. Regardless of the type of incentive, the value is the same for every respondent in a given study."
$rng = $d.Content
[void]$rng.Find.Execute("This survey is made available to respondents via", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$prng = $para.Range
$paraText = $prng.Text.TrimEnd([char]13)
[void]$prng.Find.Execute($paraText, $false, $false, $false, $false, $false, $true, 1, $false, $dynataText, 2)

# 2) "Child care" - remove surrounding gramStart/gramEnd.
FindAndClean("Child care")

# 3) "Critical retail (i.e. grocery stores, hardware stores, mechanics)" -
#    remove gramStart/gramEnd around "i.e.".
FindAndClean("Critical retail")

# 4)+5) Both "Strongly Disagree    Disagree    No Opinion    Agree    Strongly
#    Agree" paragraphs - remove spellStart/spellEnd around "Disagree".
$searchStart = 0
for ($i = 0; $i -lt 2; $i++) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute("Strongly Disagree", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { throw "Could not find occurrence $i of Strongly Disagree" }
    $para = $rng.Paragraphs(1)
    CleanParagraphProofErr $para
    $searchStart = $para.Range.End
}

# 6)-11) Education degree bullet items - remove gramStart/gramEnd around "e.g.".
FindAndClean("High school degree or equivalent")
FindAndClean("Associate degree")
FindAndClean("Bachelor")
FindAndClean("Master's degree")
FindAndClean("Professional degree")
FindAndClean("Doctorate")

# 12) "I am a full time student" - remove gramStart/gramEnd around "full time".
FindAndClean("I am a")

# 13) "Stand alone home" - remove spellStart/spellEnd around "Stand alone".
FindAndClean("Stand alone")

# ---------------------------------------------------------------------------
# 14) Add a new "Prefer not to say" list item after "Stand alone home" in the
#     housing-type list (numId 16).
# ---------------------------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute("Stand alone home", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.InsertParagraphAfter()
$idx = $para.Index
$newPara = $d.Paragraphs($idx + 1)
$newPara.Range.InsertBefore("Prefer not to say")

# ---------------------------------------------------------------------------
# 15) Add a new "Prefer not to say" list item after "Greater than $200,000" in
#     the household-income list (numId 12).
# ---------------------------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute("Greater than `$200,000", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.InsertParagraphAfter()
$idx = $para.Index
$newPara = $d.Paragraphs($idx + 1)
$newPara.Range.InsertBefore("Prefer not to say")

Write-Host "Done"
